$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete now-obsolete rows 8:10 (sending cluster "ECs" rows removed entirely)
$ws.Range("A8:T10").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Dlk1"
$ws.Range("C2").Value = "Notch1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.101448
$ws.Range("H2").Value = 6.304344
$ws.Range("I2").Value = 0.5480341737688159
$ws.Range("J2").Value = 0.5480341737688159
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 81.01644527581333
$ws.Range("R2").Value = 729.14800748232
$ws.Range("S2").Value = 0.3155663923527879
$ws.Range("T2").Value = 0.3155663923527879

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Dlk1"
$ws.Range("C3").Value = "Notch1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.101448
$ws.Range("H3").Value = 6.304344
$ws.Range("I3").Value = 0.5480341737688159
$ws.Range("J3").Value = 0.5480341737688159
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.370676
$ws.Range("N3").Value = 16.112028
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("Q3").Value = 11.286196338848
$ws.Range("R3").Value = 101.575767049632
$ws.Range("S3").Value = 0.0439607569785436
$ws.Range("T3").Value = 0.04396075697854361

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Dlk1"
$ws.Range("C4").Value = "Notch1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.101448
$ws.Range("H4").Value = 6.304344
$ws.Range("I4").Value = 0.5480341737688159
$ws.Range("J4").Value = 0.5480341737688159
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 48.39605673969333
$ws.Range("R4").Value = 435.56451065724
$ws.Range("S4").Value = 0.1885070244374843
$ws.Range("T4").Value = 0.1885070244374843

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Dlk1"
$ws.Range("C5").Value = "Notch1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.733072
$ws.Range("H5").Value = 5.199216
$ws.Range("I5").Value = 0.4519658262311841
$ws.Range("J5").Value = 0.4519658262311841
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 38.55267666666666
$ws.Range("N5").Value = 115.65803
$ws.Range("O5").Value = 0.5758151725879548
$ws.Range("P5").Value = 0.5758151725879548
$ws.Range("Q5").Value = 66.81456445605332
$ws.Range("R5").Value = 601.33108010448
$ws.Range("S5").Value = 0.2602487802351668
$ws.Range("T5").Value = 0.2602487802351668

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Dlk1"
$ws.Range("C6").Value = "Notch1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.733072
$ws.Range("H6").Value = 5.199216
$ws.Range("I6").Value = 0.4519658262311841
$ws.Range("J6").Value = 0.4519658262311841
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.370676
$ws.Range("N6").Value = 16.112028
$ws.Range("O6").Value = 0.08021535714867321
$ws.Range("P6").Value = 0.08021535714867323
$ws.Range("Q6").Value = 9.307768196672
$ws.Range("R6").Value = 83.769913770048
$ws.Range("S6").Value = 0.03625460017012961
$ws.Range("T6").Value = 0.03625460017012962

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Dlk1"
$ws.Range("C7").Value = "Notch1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.733072
$ws.Range("H7").Value = 5.199216
$ws.Range("I7").Value = 0.4519658262311841
$ws.Range("J7").Value = 0.4519658262311841
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.02986166666667
$ws.Range("N7").Value = 69.089585
$ws.Range("O7").Value = 0.3439694702633719
$ws.Range("P7").Value = 0.3439694702633719
$ws.Range("Q7").Value = 39.91240841837333
$ws.Range("R7").Value = 359.21167576536
$ws.Range("S7").Value = 0.1554624458258876
$ws.Range("T7").Value = 0.1554624458258876
